$wb = $excel.ActiveWorkbook

# Sheet "Dados": update E2:E5 text (the corrected-vs-uncorrected comparison text),
# and move the active selection to E8.
$ws = $wb.Worksheets.Item("Dados")
$newText = "o desvio padrão corrigido é maior, ou bastante maior, que o desvio padrão não corrigido"
$ws.Range("E2").Value = $newText
$ws.Range("E3").Value = $newText
$ws.Range("E4").Value = $newText
$ws.Range("E5").Value = $newText

$ws.Activate()
$ws.Range("E8").Select()
